# Preparation for publication 0.2.0
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update Version value (row 3, column B): 0.1.1 -> 0.2.0
$ws.Range("B3").Value = "0.2.0"

# 2. Update Date value (row 8, column B): 2023-10-19T16:17:18+00:00 -> 2023-10-19T17:05:12+00:00
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# 3. Insert a new row after "Contact" (row 10) for "Jurisdiction" / "iso:code:3166:FR",
#    pushing all following rows down by one.
$ws.Rows.Item(11).Insert()

# Fix up the formatting of the newly inserted row so it matches the other body rows
# (the insert operation leaves a blank/default style instead of reusing the body style).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
